$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "14-10-2021"
$ws.Range("B32").Value = 10000
$ws.Range("D32").Value = 0
